$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "age" header column
$ws.Range("E1").Value = 'age'

# Hide rows 26:40 (previously visible, now hidden like the rest of the data rows)
$ws.Rows("26:40").Hidden = $true

# Use an existing fully-styled row as the formatting template for new rows
$ws.Range("A1:E1").Copy()

# Append new data rows 41:49
# Row 41
$ws.Range("A41:E41").PasteSpecial(-4122)
$ws.Range("A41").Value = 'Steve''s Real Food Freeze-Dried Raw Nuggets - Beef Recipe - Raw Pet Food for Dogs and Cats - 20 oz'
$ws.Range("B41").Value = 'Ground beef, beef liver, beef kidney, broccoli, beef bone, carrots, apples, romaine lettuce, goat''s milk, coconut oil, sesame seeds, salmon oil, flaxseed, dried kelp, taurine, green lipped mussel, inulin, eggshell
'
$ws.Range("C41").Value = 'CONVENIENT RAW FEEDING - The easiest way to feed raw as a meal, irresistible food topper, or high-quality treat.
COMPLETE NUTRITION - Complete and balanced raw diet for both cats and dogs of all life stages.
GREAT FOR PICKY EATERS - Extremely palatable and works as a meal or topper.
PURPOSEFUL INGREDIENTS - No synthetic vitamin packs, no fillers, no artificial colors or preservatives.
THE ORIGINAL RAW PET FOOD COMPANY: The nation''s first commercially available raw pet food diet trusted since 1998.'
$ws.Range("D41").Value = 41.67
$ws.Range("E41").Value = 'Adult'
$ws.Rows("41").EntireRow.AutoFit()

# Row 42
$ws.Range("A42:E42").PasteSpecial(-4122)
$ws.Range("A42").Value = 'Wellness Appetizing Entrées Shredded Adult Wet Cat Food, Natural, Protein-Rich, Grain Free, 2.8 Ounce Pouch, 8 Pack (Chicken Recipe)'
$ws.Range("B42").Value = 'Chicken, Chicken Broth, Water Sufficient for Processing, Sunflower Oil, Salmon Oil, Pumpkin, Taurine, Calcium Chloride, Monosodium Phosphate, Dried Chicory Root, Salt, Magnesium Sulfate, Choline Chloride, Potassium Chloride, Zinc Proteinate, Vitamin E Supplement, Iron Proteinate, Thiamine Mononitrate, Niacin, Manganese Proteinate, Copper Proteinate, Vitamin A Supplement, Pyridoxine Hydrochloride, Sodium Selenite, d-Calcium Pantothenate, Riboflavin Supplement, Folic Acid, Potassium Iodide, Biotin, Vitamin B12 Supplement, Vitamin D3 Supplement.'
$ws.Range("C42").Value = 'MOUTHWATERING MEALS CATS CRAVE: Delicious, protein-rich meals featuring real chicken as the #1 ingredient now in 2.8 oz pouches. More bites for more purrs!
COMPLETE & BALANCED MEAL: Wellness Appetizing Entrées offers complete and balanced nutrition for adult cats and is crafted with real, high-quality, natural* ingredients. *with added vitamins, minerals & taurine
SUPPORTS HEALTHY HYDRATION: Moisture-rich meals made with delicious broth to support cat’s daily hydration
VARIETY OF FORMS & FLAVORS: Mix & match between mouth-watering recipes and different forms to support your cat’s daily nutrition
VETERANIARIAN RECOMMENDED**: Each recipe is crafted to support the 5 Signs of Wellbeing: Sustained Energy, Highly Digestible, Healthy Skin & Coat, Immune Health, and Healthy Teeth & Bones **Survey of US Veterinarians. Data on File.'
$ws.Range("D42").Value = 15.92
$ws.Range("E42").Value = 'Adult'
$ws.Rows("42").EntireRow.AutoFit()

# Row 43
$ws.Range("A43:E43").PasteSpecial(-4122)
$ws.Range("A43").Value = 'Freeze Dried Raw Cat and Dog Treats, Healthy Limited Ingredient Chicken Cube Weight Control Traing Treats for Small Dogs Cats 5 Ounce'
$ws.Range("B43").Value = 'Chicken'
$ws.Range("C43").Value = 'PURE NATURAL LIMITED INGREDIENT. Made by 100% premium chicken meet strict USDA standards fit for human consumption, it''s free of any corn and soy.
NUTRITIOUS & DELICIOUS. The freeze-drying technology ensures that only 95% of the moisture is removed, preserving the complete nutritional content of the chicken.
MULTIPLE FEEDING METHODS. Freeze dried can be fed directly in training, fed with water, and fed with a regular diet. These three methods can help your pet consume more other nutrients.
EASY TO CARRY. The lid is easy to close and reserve, it can keep away from moisture with double sealing and easy to take anywhere you want to take your pet to feed.
SUITABLE FOR ALL AGES DOGS AND CATS. Freeze-dried chicken cube is free of any indigestion and can be easily digested in the intestines and intestines of even elderly dogs and cats.'
$ws.Range("D43").Value = 15.09
$ws.Range("E43").Value = 'All Life Stages'
$ws.Rows("43").EntireRow.AutoFit()

# Row 44
$ws.Range("A44:E44").PasteSpecial(-4122)
$ws.Range("A44").Value = 'Ketogenic Pet Foods - Keto-Kibble – High Protein, Low Carb, Starch Free, Grain Free Dog & Cat Food - 18 lb'
$ws.Range("B44").Value = 'Chicken Meal, Organic Chicken, Chicken Fat (preserved with Mixed Tocopherols), Dried Plain Beet Pulp, Meat Protein Isolate, Gelatin, Natural Flavor, Dried Apple Fiber, Dried Natural Meat Broth, Flaxseeds, Coconut Oil, Avocado Oil, Rice Bran Oil, Taurine, Calcium Carbonate, Calcium Propionate (a preservative), Choline Chloride, Fish Oil, Yeast Extract, Citric Acid (a preservative), Inulin, Mixed Tocopherols (a preservative), Rosemary Extract, Yeast Culture, Minerals (Potassium Chloride, Zinc Proteinate, Zinc Sulfate, Ferrous Sulfate, Iron Proteinate, Copper Sulfate, Copper Proteinate, Manganese Sulfate, Manganese Proteinate, Sodium Selenite, Calcium Iodate), Vitamins (Ascorbic Acid [source of Vitamin C], Vitamin E Supplement, Niacin Supplement, Calcium Pantothenate, Vitamin A Supplement, Thiamine Mononitrate, Pyridoxine Hydrochloride, Riboflavin Supplement, Vitamin D3 Supplement, Biotin, Vitamin B12 Supplement, Folic Acid), Dried Aspergillus oryzae Fermentation Extract, Dried Enterococcus faecium Fermentation Product, Dried Lactobacillus casei Fermentation Product, Dried Lactobacillus acidophilus Fermentation Product, Dried Bacillus subtilis Fermentation Extract, Dried Lactobacillus plantarum Fermentation Product, Dried Lactobacillus lactis Fermentation Product.'
$ws.Range("C44").Value = 'Protein & fat based – Our pets are carnivores by nature. Kato-kibble is formulated to feed the carnivorous creature within! With more than twice the protein Content (60+%) of most pet Foods, keta-kibble provides a protein, fat, and micronutrient rich meal designed to encourage proper metabolic balance. Ketogenic Pet Foods are the answer to carbohydrate-based pet Foods
Starch free & Grain Free – a true alternative to carbohydrate-based pet Foods – both conventional pet Foods and those that claim to be "Grain free, " but still have high starch Content. Kato-kibble is free of inflammatory, starch fillers which can spike Insulin, disrupt metabolism, and ultimately compromise your pet''s wellbeing
Formulated for both cats & Dogs – keta-kibble is a meat-based, starch free pet Food Designed for the health of dogs & CATS. Wild canines and felines consume similarly carnivorous diets. Properly designed pet Foods like keta-kibble are equally suitable for dogs and cats
Complete meal or supplement – introduce as a supplement (add-on) and protein booster to your pet''s current Food, or use as a primary meal. (Be sure to introduce keta-kibble gradually, and follow the feeding guidelines. )
High standards - manufactured in the USA. Formulated to far exceed the minimum nutritional levels established by the AAFCO dog Food and cat Food nutrient Profiles for maintenance.'
$ws.Range("D44").Value = 92.68
$ws.Range("E44").Value = 'All Life Stages'
$ws.Rows("44").EntireRow.AutoFit()

# Row 45
$ws.Range("A45:E45").PasteSpecial(-4122)
$ws.Range("A45").Value = 'The Honest Kitchen Human Grade Whole Food Clusters Grain Free Chicken Dry Cat Food, 10 lb Bag'
$ws.Range("B45").Value = 'Chicken, peas, eggs, chicken liver, lentils, potatoes, tricalcium phosphate, natural chicken flavor, flaxseed, salmon oil, sodium chloride, vitamins [niacin supplement, vitamin e supplement, vitamin a supplement, folic acid, thiamine mononitrate (vitamin b1), calcium pantothenate (vitamin b5), biotin, pyridoxine hydrochloride (vitamin b6), vitamin d3 supplement, riboflavin (vitamin b2), vitamin b12 supplement], fenugreek seed, taurine, choline chloride, dried organic kelp, minerals [zinc amino acid chelate, iron amino acid chelate, copper amino acid chelate, manganese amino acid chelate, sodium selenite], pumpkin, blueberries, cranberries, carrots, organic barley grass, yucca schidigera extract, potassium chloride, mixed tocopherols, rosemary extract, turmeric, l-carnitine, dried bacillus coagulans fermentation product* *Contains a source of live, naturally occurring microorganisms.'
$ws.Range("C45").Value = 'HUMAN GRADE: As the first human grade dry cat food, our recipe features cage free chicken, and beneficial probiotics, providing a complete and balanced meal that supports your cat''s overall wellness
GOURMET INGREDIENTS YOU CAN TRUST: Crafted from a 100% human grade recipe with no fillers, our Whole Food Clusters dry cat food ensures your pup enjoys a nutritious meal made in a facility that meets stringent safety standards
CONVENIENT AND EASY TO DIGEST: Dry cat food clusters are crafted with human grade muscle meat & liver, oats, beneficial fats & live probiotics, making them not only a delicious choice but also easy on your cat''s digestive system, perfect for adult cats.
NO COMPROMISE ON QUALITY: We prioritize your pet''s health by saying no to feed grade ingredients, meat meals, and artificial preservatives. Our human grade chicken cat food is crafted to provide only the best nutrition without compromise
SUITABLE FOR KITTENS & ADULT CATS: Complete & balanced nutrition for adult cats and kittens of all breeds and sizes, as well as adult mothers (gestation/lactation).'
$ws.Range("D45").Value = 69.99
$ws.Range("E45").Value = 'All Life Stages'
$ws.Rows("45").EntireRow.AutoFit()

# Row 46
$ws.Range("A46:E46").PasteSpecial(-4122)
$ws.Range("A46").Value = 'Nutri Bites Freeze Dried Chicken Dog & Cat Treats | Healthy Pet Training Treats or Food Topper | All Natural, 1 Single Animal Protein, High Protein | Premium Bulk Value Pack, 15 oz'
$ws.Range("B46").Value = 'See Label'
$ws.Range("C46").Value = '100% NATURAL – Our freeze-dried dog and cat treats are all-natural and made with only 1 single animal protein, real chicken. We provide a healthy treat you and your furry friend can trust that is grain free with no additives, preservatives, or fillers.
DELICIOUS TASTE FOR TRAINING OR MEAL TOPPING – Elevate your pet''s training sessions or meals with our irresistible Freeze-Dried Treats. Bursting with delicious flavor, these treats, made from premium chicken, are a perfect topping to enhance any meal. Their savory taste ensures your puppy or cat will be motivated and eager to learn with every rewarding bite.
SUPPORT FOR YOUR PET’S HEALTH – Our high-protein, easy-to-digest freeze-dried chicken snacks are not just tasty – they''re also packed with vital benefits. Nutri Bites are packed with essential fatty acids, vitamins and minerals, including Omega-6, to support your pet’s joints, teeth, bones, organs and immune system, along with protecting hair and skin.
LESS DUST & CRUMBS – Our premium standard freeze-drying process ensures more meat and less unnecessary dust and crumbs. Give your pet the pure joy of wholesome, meat-packed goodness with every treat.
15 OZ BULK VALUE POUCH – Offering the best bulk value, each large resealable pouch is packed with 425 grams of premium treats that ensure your furry friend stays happy and healthy without breaking the bank. Treat your pet to quality nutrition at an affordable price with every pouch. Try all pet treats including; freeze-dried beef liver, salmon, and chicken!'
$ws.Range("D46").Value = 22.99
$ws.Range("E46").Value = 'All Life Stages'
$ws.Rows("46").EntireRow.AutoFit()

# Row 47
$ws.Range("A47:E47").PasteSpecial(-4122)
$ws.Range("A47").Value = 'Dr. Marty Nature''s Feast Essential Wellness Poultry Dry Cat Food 12 oz'
$ws.Range("B47").Value = 'See Label'
$ws.Range("C47").Value = 'Premium Freeze-Dried Raw Nutrition: Our cat food formula is 100% natural, packed with feline-friendly fruits, veggies, and mixed proteins to support energy, good digestion, and smooth skin with a shinier coat
High Protein, Low Carb: With 37% crude protein and only 3% crude fiber, this cat food is perfect for maintaining a healthy weight and supporting strong muscles
Omega-3 Rich: Contains 5% omega-3 fatty acids to support your cat''s heart health, brain function, and a lustrous coat
Suitable for All Life Stages: Formulated to meet the nutritional needs of cats of all ages, from kittens to seniors'
$ws.Rows("47").EntireRow.AutoFit()

# Row 48
$ws.Range("A48:E48").PasteSpecial(-4122)
$ws.Range("A48").Value = 'Wysong Epigen Canine/Feline Dry Diet - Dog/Cat Food- 5 Pound Bag (WDCFE5)'
$ws.Range("B48").Value = 'Organic Chicken, Chicken Meal, Turkey Meal, Potato Protein, Meat Protein Isolate, Chicken Fat (preserved with Mixed Tocopherols), Gelatin, Dried Plain Beet Pulp, Natural Flavor, Coconut Oil, Chia Seeds, Salt, Taurine, Calcium Carbonate, Dried Tomato Pomace, Calcium Propionate (a preservative), Choline Chloride, Dried Kelp, Organic Barley Grass Powder, Dried Blueberry Powder, Dried Yogurt, Apple Fiber, Dried Kale, Dried Spinach, Dried Carrots, Fish Oil, Yeast Extract, Citric Acid (a preservative), Inulin, Mixed Tocopherols (a preservative), Rosemary Extract, Yeast Culture, Minerals (Potassium Chloride, Zinc Proteinate, Zinc Sulfate, Ferrous Sulfate, Iron Proteinate, Copper Sulfate, Copper Proteinate, Manganese Sulfate, Manganese Proteinate, Sodium Selenite, Calcium Iodate), Vitamins (Ascorbic Acid [source of Vitamin C], Vitamin E Supplement, Niacin Supplement, Calcium Pantothenate, Vitamin A Supplement, Thiamine Mononitrate, Pyridoxine Hydrochloride, Riboflavin Supplement, Vitamin D3 Supplement, Biotin, Vitamin B12 Supplement, Folic Acid), Dried Aspergillus oryzae Fermentation Extract, Dried Enterococcus faecium Fermentation Product, Dried Lactobacillus casei Fermentation Product, Dried Lactobacillus acidophilus Fermentation Product, Dried Bacillus subtilis Fermentation Extract, Dried Lactobacillus plantarum Fermentation Product, Dried Lactobacillus lactis Fermentation Product, Pepper.'
$ws.Range("C48").Value = 'The first extruded dry diet kibble pet food that is Starch Free Epigen does what "Grain Free" dog food and cat foods only pretend to do
More closely resembles the food your pet is designed to eat Has the natural flavor and taste that pets truly desire
Is replete with all the critical food elements Wysong was first to put into pet foods probiotics prebiotics enzymes phytonutrients nutraceuticals natural form antioxidants vitamins and minerals
Contains unprecedented and unequalled meat content including organic Rich in omega 3''s and other essential fatty acids
Very nutrient dense Introduce slowly according to the feeding guidelines Can be used as a protein booster and meal supplement
Manufactured by Wysong in accordance with strict quality control protocols
Wysong has been a leader in pet nutrition since 1979'
$ws.Range("D48").Value = 19.09
$ws.Range("E48").Value = 'Adult'
$ws.Rows("48").EntireRow.AutoFit()

# Row 49
$ws.Range("A49:E49").PasteSpecial(-4122)
$ws.Range("A49").Value = 'Stella & Chewy''s – Stella’s Solutions Digestive Boost – Cage-Free Chicken Dinner Mixer – Freeze-Dried Raw, Protein Rich, Grain Free Cat Food – 7.5 oz Bag'
$ws.Range("B49").Value = 'Chicken with bone, chicken liver, pumpkin, chicken heart, chicken gizzard, salmon oil, salt, choline chloride, dandelion, tocopherols (preservative), papaya, taurine, pineapple, dried chicory root (source of inulin), L-Carnitine, dried Pediococcus acidilactici fermentation product, dried Lactobacillus acidophilus fermentation product, dried Bifidobacterium longum fermentation product, dried Bacillus coagulans fermentation product, potassium chloride, sodium phosphate, dried kelp, zinc proteinate, iron proteinate, vitamin A supplement, vitamin E supplement, niacin supplement, copper proteinate, riboflavin supplement, sodium selenite, d-calcium pantothenate, biotin, manganese proteinate, thiamine mononitrate, pyridoxine hydrochloride, vitamin D3 supplement, folic acid, vitamin B12 supplement.'
$ws.Range("C49").Value = 'Ultimate Digestive Support: This recipe is chock full of ingredients to support your kitty’s digestion and intestinal health, including L-carnitine, a chemical that improves metabolism and other body processes by helping turn fat into energy
All the Fiber Your Cat Needs: Natural sources of soluble fiber, such as pumpkin and pineapple stems, help your cat digest food and treats effectively
Full of Probiotics, Prebiotics & Digestive Enzymes: We’ve added probiotics and prebiotics to aid in intestinal health and optimal digestion; papaya, a natural source of papain, also has enzymes that aid in digestion
Picky Eater Approved: Even the pickiest eaters go wild for this topper; these special and unique recipes add a boost of nutrition, raw power, and taste to your cat’s bowl!
A Healthy & Tasty Meal Mixer, Supplement or Snack: Made with cage-free chicken, this all-natural, freeze-dried raw recipe is a convenient way to add the power, nutrition and taste of raw to your cat’s diet
Supports Whole Body Health: Fortified with vitamins and minerals, our Stella Solution’s formulas help improve your cat’s entire well-being, with Omega fatty acids for skin and coat support, antioxidants for immune support, and taurine for heart health
Whole Prey Ingredients: Feeding your pet a diet similar to what their ancestors enjoyed in the wild helps your pet thrive; a raw diet helps support healthy digestion, strong teeth and gums, vibrant skin and coat, and stamina and vitality'
$ws.Range("D49").Value = 22.99
$ws.Range("E49").Value = 'Adult'
$ws.Rows("49").EntireRow.AutoFit()

